$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-log entry for 2026/02/15 (日) needs to be inserted before the
# existing row 795, shifting all subsequent rows (795-836) down by one and
# extending the used range from D836 to D837.
$ws.Rows.Item(795).EntireRow.Insert()

# Populate the newly inserted row. Column A holds the date as literal text
# (matching the rest of the sheet, which stores dates as strings rather than
# real date serials), so a leading apostrophe forces text entry instead of
# Excel's automatic date conversion.
$ws.Range("A795").Value = "'2026/02/15"
$ws.Range("B795").Value = "日"
$ws.Range("C795").Value = 0
$ws.Range("D795").Value = 22
